$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.469.22"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.30%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.631.82"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.60%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.57"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.12"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.11%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.534"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.39%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.630.62"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.64%  "

$ws.Range("E10").Value = "  -3.30%  "

$ws.Range("E11").Value = "  +1.26%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.364"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.86%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.23"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.68"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.68%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.110.07"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.67%  "

$ws.Range("E16").Value = "  -2.08%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.465.09"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.622.73"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.98"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.96%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.05"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "357.26"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.85%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.32"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.78%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.67"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.38%  "

$ws.Range("E24").Value = "  -4.60%  "

$ws.Range("E25").Value = "  +0.09%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.32"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.23%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "69.66"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.87%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.763.47"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.21%  "

$ws.Range("E30").Value = "  -1.72%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "547.86"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.73%  "

$ws.Range("E32").Value = "  -1.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.36"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.92%  "

$ws.Range("E34").Value = "  -1.96%  "

$ws.Range("E35").Value = "  +4.43%  "

$ws.Range("E36").Value = "  +0.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.50"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.90%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "156.47"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.03"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.64%  "

$ws.Range("E40").Value = "  -2.28%  "

$ws.Range("E41").Value = "  -0.64%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "18.30"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.99%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.22"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.74%  "

$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("E45").Value = "  -3.91%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0299"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "153.07"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.36%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.579"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.80"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.51%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.71"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.21%  "

$ws.Range("E51").Value = "  -1.11%  "
